$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure cells keep their literal text representation (no auto numeric coercion,
# so values like "0.8100" or "0.00006799" keep trailing/leading zeros exactly).
$cells = @("D2", "D3", "D4", "D6", "D7", "D8", "D9", "D10", "D11", "D12", "D13", "D14", "D15", "D16", "D17", "D27", "D40", "D41", "D42", "D43", "D45", "D47", "D48", "B18", "C18", "D18", "E18", "B19", "C19", "D19", "E19", "B20", "C20", "D20", "E20", "B21", "C21", "D21", "E21", "B22", "C22", "D22", "E22", "B23", "C23", "D23", "E23", "B24", "C24", "D24", "E24", "D44", "E44")
foreach ($cellRef in $cells) {
    $ws.Range($cellRef).NumberFormat = "@"
}

$ws.Range("D2").Value = "242.49"
$ws.Range("D3").Value = "22.94"
$ws.Range("D4").Value = "5.414"
$ws.Range("D6").Value = "3.438"
$ws.Range("D7").Value = "6.544"
$ws.Range("D8").Value = "0.8100"
$ws.Range("D9").Value = "0.9537"
$ws.Range("D10").Value = "0.1422"
$ws.Range("D11").Value = "0.07430"
$ws.Range("D12").Value = "0.03275"
$ws.Range("D13").Value = "0.03053"
$ws.Range("D14").Value = "0.09331"
$ws.Range("D15").Value = "3.859"
$ws.Range("D16").Value = "0.001577"
$ws.Range("D17").Value = "0.04670"
$ws.Range("D27").Value = "0.0002284"
$ws.Range("D40").Value = "0.03931"
$ws.Range("D41").Value = "0.006186"
$ws.Range("D42").Value = "0.1068"
$ws.Range("D43").Value = "0.003000"
$ws.Range("D45").Value = "0.00005197"
$ws.Range("D47").Value = "0.7199"
$ws.Range("D48").Value = "0.002392"
$ws.Range("B18").Value = "TigerCash"
$ws.Range("C18").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("D18").Value = "0.005864"
$ws.Range("E18").Value = "17TigerCashTCH"
$ws.Range("B19").Value = "BitKan"
$ws.Range("C19").Value = "https://coinranking.com/coin/RDOsLDgvY-AXe+bitkan-kan"
$ws.Range("D19").Value = "0.001259"
$ws.Range("E19").Value = "18BitKanKAN"
$ws.Range("B20").Value = "HotbitToken"
$ws.Range("C20").Value = "https://coinranking.com/coin/uQJB8Ocu8lTb+hotbittoken-htb"
$ws.Range("D20").Value = "0.004896"
$ws.Range("E20").Value = "19HotbitTokenHTB"
$ws.Range("B21").Value = "NitroEx"
$ws.Range("C21").Value = "https://coinranking.com/coin/8oiZw6gwYhC+nitroex-ntx"
$ws.Range("D21").Value = "0.00006799"
$ws.Range("E21").Value = "20NitroExNTX"
$ws.Range("B22").Value = "LEO"
$ws.Range("C22").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D22").Value = "3.574"
$ws.Range("E22").Value = "21LEOLEO"
$ws.Range("B23").Value = "BTSEToken"
$ws.Range("C23").Value = "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
$ws.Range("D23").Value = "2.132"
$ws.Range("E23").Value = "22BTSETokenBTSE"
$ws.Range("B24").Value = "One"
$ws.Range("C24").Value = "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
$ws.Range("D24").Value = "0.01129"
$ws.Range("E24").Value = "23OneONEBestin24h"
$ws.Range("D44").Value = "0.009862"
$ws.Range("E44").Value = "43LocalTradersLCT"
